$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header/count values)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (meanEMG / legmaxROM values)
$ws.Range("B2").Value = 16.150000000000002
$ws.Range("C2").Value = 11.450000000000001
$ws.Range("D2").Value = 13.45
$ws.Range("E2").Value = 13.25

# Row 3
$ws.Range("B3").Value = 11.950000000000001
$ws.Range("C3").Value = 6.8500000000000005
$ws.Range("D3").Value = 9.75
$ws.Range("E3").Value = 13.05

# Update selection to match the new reduced range of interest
$ws.Range("B1:E3").Select()
